$wb = $excel.ActiveWorkbook

$wsWeapons = $wb.Worksheets.Item("Weapons")
$wsBoxes   = $wb.Worksheets.Item("Boxes")
$wsMags    = $wb.Worksheets.Item("Mags")

# --- Weapons sheet: update Qty (column B) balance numbers ---
$wsWeapons.Range("B2").Value  = 0   # 9mm Pistol
$wsWeapons.Range("B3").Value  = 0   # 10mm Pistol
$wsWeapons.Range("B4").Value  = 2   # Hunting Rifle
$wsWeapons.Range("B5").Value  = 0   # Pipe Pistol
$wsWeapons.Range("B7").Value  = 0   # Laser Pistol
$wsWeapons.Range("B8").Value  = 2   # Heavy Laser Pistol
$wsWeapons.Range("B9").Value  = 2   # Laser Rifle
$wsWeapons.Range("B11").Value = 3   # S&W 357 Magnum
$wsWeapons.Range("B14").Value = 2   # Combat Shotgun
$wsWeapons.Range("J14").Value = 7   # Combat Shotgun, column J

# --- Boxes sheet: update Qty (column B) balance numbers ---
$wsBoxes.Range("B2").Value = 1   # 9mm
$wsBoxes.Range("B3").Value = 1   # 10mm
$wsBoxes.Range("B4").Value = 1   # x308
$wsBoxes.Range("B5").Value = 2   # x38
$wsBoxes.Range("B6").Value = 2   # x357
$wsBoxes.Range("B9").Value = 0   # 12 Gauge

# --- Mags sheet: update Qty (column C) balance numbers ---
$wsMags.Range("C2").Value = 4   # 9mm Pistol Mag
$wsMags.Range("C3").Value = 4   # 10mm Pistol Mag
$wsMags.Range("C4").Value = 2   # Hunting Rifle Mag
$wsMags.Range("C7").Value = 4   # Energy Cell
$wsMags.Range("C8").Value = 2   # Fusion Cell

# --- Sheet view / selection updates (card resize) ---
# Select on Weapons and Mags first (without leaving them as the active tab),
# then finish on Boxes so it ends up the active/selected sheet tab.
$wsWeapons.Range("B12").Select()
$wsMags.Range("C9").Select()
$wsBoxes.Range("B5").Select()
